$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing data rows (2..10) down to (3..11)
$ws.Rows("2:2").Insert()

# The inserted row picks up header-row formatting by default; clear it so the
# new data row matches the plain (unstyled) look of the other data rows.
$ws.Range("A2:T2").ClearFormats()

# Re-apply the date number format used by the other rows' "Fecha" column (D).
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row with the new price record.
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C2").Value = "Los Lagos"
$ws.Range("D2").Value = 44530
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 500
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("R2").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S2").Value = 1139
$ws.Range("T2").Value = 18
